$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in this sheet, even when they
# look like plain numbers (e.g. "0.9996"). A leading apostrophe forces Excel
# to keep the assigned value as text instead of auto-converting it to a
# number; the apostrophe itself is not stored. Values that already contain a
# second "." (e.g. "29.165.94") are never number-like and need no prefix.
$ws.Range("D2").Value = '29.165.94'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.841.19'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").Value = '''0.9996'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''241.09'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").Value = '''0.6862'
$ws.Range("E6").Value = '  -2.49%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '''0.2997'
$ws.Range("E8").Value = '  -2.25%  '
$ws.Range("E9").Value = '  -3.29%  '
$ws.Range("D10").Value = '''23.19'
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").Value = '''0.07654'
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("D12").Value = '1.845.18'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").Value = '''5.053'
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").Value = '''0.6815'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").Value = '''87.21'
$ws.Range("E15").Value = '  -6.57%  '
$ws.Range("D16").Value = '''6.151'
$ws.Range("E16").Value = '  -6.85%  '
$ws.Range("D17").Value = '29.164.77'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").Value = '''0.000008186'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '2.084.23'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '''228.80'
$ws.Range("E20").Value = '  -5.53%  '
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("D22").Value = '''0.9997'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '''7.384'
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''159.67'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '''0.1445'
$ws.Range("E26").Value = '  -3.98%  '
$ws.Range("D27").Value = '''8.753'
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").Value = '''18.07'
$ws.Range("E28").Value = '  -1.25%  '
$ws.Range("D29").Value = '''1.512'
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("D30").Value = '''4.276'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '''0.05250'
$ws.Range("E33").Value = '  +2.41%  '
$ws.Range("D34").Value = '''0.7593'
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("D35").Value = '''1.852'
$ws.Range("E35").Value = '  -2.60%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '1.300.90'
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("D39").Value = '''0.01831'
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("D40").Value = '''2.721'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '''0.9373'
$ws.Range("E41").Value = '  -2.19%  '
$ws.Range("D42").Value = '''5.956'
$ws.Range("E42").Value = '  -1.56%  '
$ws.Range("D43").Value = '''104.89'
$ws.Range("E43").Value = '  -1.91%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '1.985.04'
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").Value = '''0.5197'
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("D47").Value = '''64.78'
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''9.476'
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.769'
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("B50").Value = 'XinFinNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D50").Value = '''0.07462'
$ws.Range("E50").Value = '  +18.15%  '
$ws.Range("E51").Value = '  +0.64%  '
